# Auto-generated from the target OOXML diff.
# Applies per-cell numeric updates (+ two additions, one deletion)
# to the FFXIV "Brynhildr_Profits" crafting-profit workbook, sheet by sheet.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4064.8667
$ws.Range("J64").Value = 4069.5
$ws.Range("L64").Value = 4069.5
$ws.Range("N64").Value = -4565.5
$ws.Range("H67").Value = 4064.8667
$ws.Range("J67").Value = 4069.5
$ws.Range("L67").Value = 4069.5
$ws.Range("N67").Value = -5785.5
$ws.Range("H76").Value = 3456.9636
$ws.Range("I76").Value = 3141.9167
$ws.Range("J76").Value = 4053.8948
$ws.Range("K76").Value = 3141.9167
$ws.Range("L76").Value = 4053.8948
$ws.Range("M76").Value = -2826.9167
$ws.Range("N76").Value = -4683.8948
$ws.Range("H79").Value = 3456.9636
$ws.Range("I79").Value = 3141.9167
$ws.Range("J79").Value = 4053.8948
$ws.Range("K79").Value = 3141.9167
$ws.Range("L79").Value = 4053.8948
$ws.Range("M79").Value = -2049.9167
$ws.Range("N79").Value = -6237.8948
$ws.Range("H86").Value = 6808
$ws.Range("I86").Value = 2229.5
$ws.Range("J86").Value = 13675.75
$ws.Range("K86").Value = 2229.5
$ws.Range("L86").Value = 13675.75
$ws.Range("M86").Value = -1106.5
$ws.Range("N86").Value = -15921.75
$ws.Range("H89").Value = 6808
$ws.Range("I89").Value = 2229.5
$ws.Range("J89").Value = 13675.75
$ws.Range("K89").Value = 11147.5
$ws.Range("L89").Value = 68378.75
$ws.Range("M89").Value = -5531.5
$ws.Range("N89").Value = -79610.75
$ws.Range("H113").Value = 3928.111
$ws.Range("I113").Value = 3542.7144
$ws.Range("J113").Value = 5277
$ws.Range("K113").Value = 3542.7144
$ws.Range("L113").Value = 5277
$ws.Range("M113").Value = -288.7143999999998
$ws.Range("N113").Value = -11785
$ws.Range("H116").Value = 26410.5
$ws.Range("I116").Value = 13136.625
$ws.Range("K116").Value = 13136.625
$ws.Range("M116").Value = -9694.625
$ws.Range("H132").Value = 7675.2256
$ws.Range("I132").Value = 8334.286
$ws.Range("J132").Value = 1524
$ws.Range("K132").Value = 25002.858
$ws.Range("L132").Value = 4572
$ws.Range("M132").Value = -22472.858
$ws.Range("N132").Value = -9632

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2037.8695
$ws.Range("I2").Value = 1745.7222
$ws.Range("K2").Value = 1745.7222
$ws.Range("M2").Value = -1632.7222
$ws.Range("H88").Value = 1805.7368
$ws.Range("I88").Value = 1894.6666
$ws.Range("K88").Value = 1894.6666
$ws.Range("M88").Value = -1488.6666
$ws.Range("H91").Value = 1805.7368
$ws.Range("I91").Value = 1894.6666
$ws.Range("K91").Value = 1894.6666
$ws.Range("M91").Value = -490.6666
$ws.Range("H116").Value = 2037.8695
$ws.Range("I116").Value = 1745.7222
$ws.Range("K116").Value = 1745.7222
$ws.Range("M116").Value = 548.2778000000001
$ws.Range("H122").Value = 1925.0714
$ws.Range("I122").Value = 1803.9231
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 5411.7693
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -2961.7693
$ws.Range("N122").Value = -15400
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2037.8695
$ws.Range("I3").Value = 1745.7222
$ws.Range("K3").Value = 1745.7222
$ws.Range("M3").Value = -1631.7222
$ws.Range("H114").Value = 47250
$ws.Range("J114").Value = 47250
$ws.Range("L114").Value = 47250
$ws.Range("N114").Value = -55928
$ws.Range("H115").Value = 34275.6
$ws.Range("J115").Value = 34275.6
$ws.Range("L115").Value = 34275.6
$ws.Range("N115").Value = -37409.6

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2194410.8
$ws.Range("I31").Value = 3760818.5
$ws.Range("K31").Value = 3760818.5
$ws.Range("M31").Value = -3760523.5
$ws.Range("H34").Value = 2194410.8
$ws.Range("I34").Value = 3760818.5
$ws.Range("K34").Value = 3760818.5
$ws.Range("M34").Value = -3760616.5
$ws.Range("H99").Value = 10210.519
$ws.Range("I99").Value = 18229
$ws.Range("J99").Value = 2764.7856
$ws.Range("K99").Value = 18229
$ws.Range("L99").Value = 2764.7856
$ws.Range("M99").Value = -16731
$ws.Range("N99").Value = -5760.7856
$ws.Range("H122").Value = 8514.526
$ws.Range("I122").Value = 1895.7715
$ws.Range("K122").Value = 5687.3145
$ws.Range("M122").Value = -3237.3145
$ws.Range("H126").Value = 10210.519
$ws.Range("I126").Value = 18229
$ws.Range("J126").Value = 2764.7856
$ws.Range("K126").Value = 54687
$ws.Range("L126").Value = 8294.356800000001
$ws.Range("M126").Value = -52217
$ws.Range("N126").Value = -13234.3568

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1552791
$ws.Range("J5").Value = 1369197.2
$ws.Range("L5").Value = 4107591.6
$ws.Range("N5").Value = -4107815.6
$ws.Range("H37").Value = 68529.414
$ws.Range("J37").Value = 68529.414
$ws.Range("L37").Value = 205588.242
$ws.Range("N37").Value = -205812.242
$ws.Range("H68").Value = 7846.129
$ws.Range("I68").Value = 2097.6667
$ws.Range("J68").Value = 9225.76
$ws.Range("K68").Value = 6293.000100000001
$ws.Range("L68").Value = 27677.28
$ws.Range("M68").Value = -5482.000100000001
$ws.Range("N68").Value = -29299.28
$ws.Range("H71").Value = 7846.129
$ws.Range("I71").Value = 2097.6667
$ws.Range("J71").Value = 9225.76
$ws.Range("K71").Value = 18879.0003
$ws.Range("L71").Value = 83031.84
$ws.Range("M71").Value = -14823.0003
$ws.Range("N71").Value = -91143.84
$ws.Range("H75").Value = 700
$ws.Range("I75").Value = 500
$ws.Range("J75").Value = 900
$ws.Range("K75").Value = 1500
$ws.Range("L75").Value = 2700
$ws.Range("M75").Value = -502
$ws.Range("N75").Value = -4696
$ws.Range("H78").Value = 700
$ws.Range("I78").Value = 500
$ws.Range("J78").Value = 900
$ws.Range("K78").Value = 4500
$ws.Range("L78").Value = 8100
$ws.Range("M78").Value = 492
$ws.Range("N78").Value = -18084
$ws.Range("H86").Value = 789.625
$ws.Range("I86").Value = 780.6667
$ws.Range("K86").Value = 2342.0001
$ws.Range("M86").Value = -1156.0001
$ws.Range("H89").Value = 789.625
$ws.Range("I89").Value = 780.6667
$ws.Range("K89").Value = 7026.0003
$ws.Range("M89").Value = -1098.0003
$ws.Range("H122").Value = 1009125.56
$ws.Range("I122").Value = 3226001.5
$ws.Range("J122").Value = 1454.7273
$ws.Range("K122").Value = 29034013.5
$ws.Range("L122").Value = 13092.5457
$ws.Range("M122").Value = -29031563.5
$ws.Range("N122").Value = -17992.5457
$ws.Range("H135").Value = 1552791
$ws.Range("J135").Value = 1369197.2
$ws.Range("L135").Value = 12322774.8
$ws.Range("N135").Value = -12327844.8
$ws.Range("H137").Value = 6350.2
$ws.Range("I137").Value = 1843
$ws.Range("J137").Value = 11859
$ws.Range("K137").Value = 5529
$ws.Range("L137").Value = 35577
$ws.Range("M137").Value = -429
$ws.Range("N137").Value = -45777

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1264
$ws.Range("I97").Value = 1326.2354
$ws.Range("J97").Value = 999.5
$ws.Range("K97").Value = 1326.2354
$ws.Range("L97").Value = 999.5
$ws.Range("M97").Value = -830.2354
$ws.Range("N97").Value = -1991.5
$ws.Range("H102").Value = 2562.8518
$ws.Range("I102").Value = 2624.739
$ws.Range("K102").Value = 2624.739
$ws.Range("M102").Value = -1002.739
$ws.Range("H122").Value = 4515.6924
$ws.Range("I122").Value = 5764.2856
$ws.Range("J122").Value = 3059
$ws.Range("K122").Value = 17292.8568
$ws.Range("L122").Value = 9177
$ws.Range("M122").Value = -14842.8568
$ws.Range("N122").Value = -14077

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 424.66666
$ws.Range("I113").Value = 573
$ws.Range("K113").Value = 1719
$ws.Range("M113").Value = 451

Write-Host "Applied 203 cell updates across 7 sheets."
